$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 117-120 (re-ordering of 4 matches) ---
# Row 117
$ws.Cells.Item(117, 2).Value = 7013886
$ws.Cells.Item(117, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(117, 4).Value = "Uruguay Clausura"
$ws.Cells.Item(117, 5).Value = 45267.70833333334
$ws.Cells.Item(117, 6).Value = "Racing Club de Montevideo"
$ws.Cells.Item(117, 7).Value = "Cerro"
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 9).Value = 1
$ws.Cells.Item(117, 10).Value = "A"
$ws.Cells.Item(117, 11).Value = 2.25
$ws.Cells.Item(117, 12).Value = 3.1
$ws.Cells.Item(117, 13).Value = 3.25
$ws.Cells.Item(117, 14).Value = 2.25
$ws.Cells.Item(117, 15).Value = 2.875
$ws.Cells.Item(117, 16).Value = 3.5
$ws.Cells.Item(117, 17).Value = -0.25
$ws.Cells.Item(117, 18).Value = 1.95
$ws.Cells.Item(117, 19).Value = 1.9
$ws.Cells.Item(117, 20).Value = 2
$ws.Cells.Item(117, 21).Value = 1.925
$ws.Cells.Item(117, 22).Value = 1.925
$ws.Cells.Item(117, 23).Value = -1
$ws.Cells.Item(117, 24).Value = -1
$ws.Cells.Item(117, 25).Value = 2.5
$ws.Cells.Item(117, 26).Value = -1
$ws.Cells.Item(117, 27).Value = 0.8999999999999999
$ws.Cells.Item(117, 28).Value = -1
$ws.Cells.Item(117, 29).Value = 0.925

# Row 118
$ws.Cells.Item(118, 2).Value = 7013409
$ws.Cells.Item(118, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(118, 4).Value = "Uruguay Clausura"
$ws.Cells.Item(118, 5).Value = 45267.70833333334
$ws.Cells.Item(118, 6).Value = "Nacional De Football"
$ws.Cells.Item(118, 7).Value = "Torque"
$ws.Cells.Item(118, 8).Value = 1
$ws.Cells.Item(118, 9).Value = 1
$ws.Cells.Item(118, 10).Value = "D"
$ws.Cells.Item(118, 11).Value = 1.666
$ws.Cells.Item(118, 12).Value = 3.9
$ws.Cells.Item(118, 13).Value = 4.5
$ws.Cells.Item(118, 14).Value = 1.615
$ws.Cells.Item(118, 15).Value = 4
$ws.Cells.Item(118, 16).Value = 4.75
$ws.Cells.Item(118, 17).Value = -0.75
$ws.Cells.Item(118, 18).Value = 1.8
$ws.Cells.Item(118, 19).Value = 2.05
$ws.Cells.Item(118, 20).Value = 2.75
$ws.Cells.Item(118, 21).Value = 1.95
$ws.Cells.Item(118, 22).Value = 1.9
$ws.Cells.Item(118, 23).Value = -1
$ws.Cells.Item(118, 24).Value = 3
$ws.Cells.Item(118, 25).Value = -1
$ws.Cells.Item(118, 26).Value = -1
$ws.Cells.Item(118, 27).Value = 1.05
$ws.Cells.Item(118, 28).Value = -1
$ws.Cells.Item(118, 29).Value = 0.8999999999999999

# Row 119
$ws.Cells.Item(119, 2).Value = 7013885
$ws.Cells.Item(119, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(119, 4).Value = "Uruguay Clausura"
$ws.Cells.Item(119, 5).Value = 45267.70833333334
$ws.Cells.Item(119, 6).Value = "La Luz"
$ws.Cells.Item(119, 7).Value = "Atletico Fenix Montevideo"
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 9).Value = 2
$ws.Cells.Item(119, 10).Value = "A"
$ws.Cells.Item(119, 11).Value = 3
$ws.Cells.Item(119, 12).Value = 3
$ws.Cells.Item(119, 13).Value = 2.4
$ws.Cells.Item(119, 14).Value = 2.9
$ws.Cells.Item(119, 15).Value = 2.75
$ws.Cells.Item(119, 16).Value = 2.6
$ws.Cells.Item(119, 17).Value = 0
$ws.Cells.Item(119, 18).Value = 2.025
$ws.Cells.Item(119, 19).Value = 1.825
$ws.Cells.Item(119, 20).Value = 2
$ws.Cells.Item(119, 21).Value = 2.025
$ws.Cells.Item(119, 22).Value = 1.825
$ws.Cells.Item(119, 23).Value = -1
$ws.Cells.Item(119, 24).Value = -1
$ws.Cells.Item(119, 25).Value = 1.6
$ws.Cells.Item(119, 26).Value = -1
$ws.Cells.Item(119, 27).Value = 0.825
$ws.Cells.Item(119, 28).Value = 0
$ws.Cells.Item(119, 29).Value = -0

# Row 120
$ws.Cells.Item(120, 2).Value = 7013702
$ws.Cells.Item(120, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(120, 4).Value = "Uruguay Clausura"
$ws.Cells.Item(120, 5).Value = 45267.70833333334
$ws.Cells.Item(120, 6).Value = "Defensor Sporting"
$ws.Cells.Item(120, 7).Value = "Danubio"
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 9).Value = 2
$ws.Cells.Item(120, 10).Value = "A"
$ws.Cells.Item(120, 11).Value = 1.8
$ws.Cells.Item(120, 12).Value = 3.6
$ws.Cells.Item(120, 13).Value = 4.2
$ws.Cells.Item(120, 14).Value = 1.8
$ws.Cells.Item(120, 15).Value = 3.6
$ws.Cells.Item(120, 16).Value = 4.2
$ws.Cells.Item(120, 17).Value = -0.75
$ws.Cells.Item(120, 18).Value = 2.05
$ws.Cells.Item(120, 19).Value = 1.8
$ws.Cells.Item(120, 20).Value = 2.25
$ws.Cells.Item(120, 21).Value = 1.85
$ws.Cells.Item(120, 22).Value = 2
$ws.Cells.Item(120, 23).Value = -1
$ws.Cells.Item(120, 24).Value = -1
$ws.Cells.Item(120, 25).Value = 3.2
$ws.Cells.Item(120, 26).Value = -1
$ws.Cells.Item(120, 27).Value = 0.8
$ws.Cells.Item(120, 28).Value = -0.5
$ws.Cells.Item(120, 29).Value = 0.5

# --- Update rows 155-158 (existing matches updated with results) ---
# Row 155
$ws.Cells.Item(155, 1).Value = 153
$ws.Cells.Item(155, 2).Value = 7990770
$ws.Cells.Item(155, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(155, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(155, 5).Value = 45374.5625
$ws.Cells.Item(155, 6).Value = "Rampla Juniors"
$ws.Cells.Item(155, 7).Value = "Miramar Misiones"
$ws.Cells.Item(155, 8).Value = 1
$ws.Cells.Item(155, 9).Value = 1
$ws.Cells.Item(155, 10).Value = "D"
$ws.Cells.Item(155, 11).Value = 2.8
$ws.Cells.Item(155, 12).Value = 3.2
$ws.Cells.Item(155, 13).Value = 2.4
$ws.Cells.Item(155, 14).Value = 2.8
$ws.Cells.Item(155, 15).Value = 3.2
$ws.Cells.Item(155, 16).Value = 2.4
$ws.Cells.Item(155, 17).Value = 0
$ws.Cells.Item(155, 18).Value = 2.1
$ws.Cells.Item(155, 19).Value = 1.775
$ws.Cells.Item(155, 20).Value = 2.5
$ws.Cells.Item(155, 21).Value = 2.1
$ws.Cells.Item(155, 22).Value = 1.775
$ws.Cells.Item(155, 23).Value = -1
$ws.Cells.Item(155, 24).Value = 2.2
$ws.Cells.Item(155, 25).Value = -1
$ws.Cells.Item(155, 26).Value = 0
$ws.Cells.Item(155, 27).Value = -0
$ws.Cells.Item(155, 28).Value = -1
$ws.Cells.Item(155, 29).Value = 0.7749999999999999

# Row 156
$ws.Cells.Item(156, 1).Value = 154
$ws.Cells.Item(156, 2).Value = 7990771
$ws.Cells.Item(156, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(156, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(156, 5).Value = 45374.77083333334
$ws.Cells.Item(156, 6).Value = "Nacional De Football"
$ws.Cells.Item(156, 7).Value = "Club Atletico Progreso"
$ws.Cells.Item(156, 8).Value = 0
$ws.Cells.Item(156, 9).Value = 0
$ws.Cells.Item(156, 10).Value = "D"
$ws.Cells.Item(156, 11).Value = 1.4
$ws.Cells.Item(156, 12).Value = 4.5
$ws.Cells.Item(156, 13).Value = 6.5
$ws.Cells.Item(156, 14).Value = 1.333
$ws.Cells.Item(156, 15).Value = 4.5
$ws.Cells.Item(156, 16).Value = 7.5
$ws.Cells.Item(156, 17).Value = -1.25
$ws.Cells.Item(156, 18).Value = 1.9
$ws.Cells.Item(156, 19).Value = 1.95
$ws.Cells.Item(156, 20).Value = 2.5
$ws.Cells.Item(156, 21).Value = 1.95
$ws.Cells.Item(156, 22).Value = 1.9
$ws.Cells.Item(156, 23).Value = -1
$ws.Cells.Item(156, 24).Value = 3.5
$ws.Cells.Item(156, 25).Value = -1
$ws.Cells.Item(156, 26).Value = -1
$ws.Cells.Item(156, 27).Value = 0.95
$ws.Cells.Item(156, 28).Value = -1
$ws.Cells.Item(156, 29).Value = 0.8999999999999999

# Row 157
$ws.Cells.Item(157, 1).Value = 155
$ws.Cells.Item(157, 2).Value = 7990728
$ws.Cells.Item(157, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(157, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(157, 5).Value = 45374.875
$ws.Cells.Item(157, 6).Value = "Deportivo Maldonado"
$ws.Cells.Item(157, 7).Value = "Defensor Sporting"
$ws.Cells.Item(157, 8).Value = 1
$ws.Cells.Item(157, 9).Value = 0
$ws.Cells.Item(157, 10).Value = "H"
$ws.Cells.Item(157, 11).Value = 4.2
$ws.Cells.Item(157, 12).Value = 3.3
$ws.Cells.Item(157, 13).Value = 1.909
$ws.Cells.Item(157, 14).Value = 3.8
$ws.Cells.Item(157, 15).Value = 3.4
$ws.Cells.Item(157, 16).Value = 1.95
$ws.Cells.Item(157, 17).Value = 0.5
$ws.Cells.Item(157, 18).Value = 1.85
$ws.Cells.Item(157, 19).Value = 2
$ws.Cells.Item(157, 20).Value = 2.25
$ws.Cells.Item(157, 21).Value = 1.8
$ws.Cells.Item(157, 22).Value = 2.05
$ws.Cells.Item(157, 23).Value = 2.8
$ws.Cells.Item(157, 24).Value = -1
$ws.Cells.Item(157, 25).Value = -1
$ws.Cells.Item(157, 26).Value = 0.8500000000000001
$ws.Cells.Item(157, 27).Value = -1
$ws.Cells.Item(157, 28).Value = -1
$ws.Cells.Item(157, 29).Value = 1.05

# Row 158
$ws.Cells.Item(158, 1).Value = 156
$ws.Cells.Item(158, 2).Value = 7990772
$ws.Cells.Item(158, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(158, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(158, 5).Value = 45375.41666666666
$ws.Cells.Item(158, 6).Value = "Cerro"
$ws.Cells.Item(158, 7).Value = "Cerro Largo"
$ws.Cells.Item(158, 8).Value = 1
$ws.Cells.Item(158, 9).Value = 1
$ws.Cells.Item(158, 10).Value = "D"
$ws.Cells.Item(158, 11).Value = 2.75
$ws.Cells.Item(158, 12).Value = 3
$ws.Cells.Item(158, 13).Value = 2.625
$ws.Cells.Item(158, 14).Value = 2.875
$ws.Cells.Item(158, 15).Value = 3
$ws.Cells.Item(158, 16).Value = 2.55
$ws.Cells.Item(158, 17).Value = 0
$ws.Cells.Item(158, 18).Value = 2.025
$ws.Cells.Item(158, 19).Value = 1.825
$ws.Cells.Item(158, 20).Value = 2
$ws.Cells.Item(158, 21).Value = 2.05
$ws.Cells.Item(158, 22).Value = 1.8
$ws.Cells.Item(158, 23).Value = -1
$ws.Cells.Item(158, 24).Value = 2
$ws.Cells.Item(158, 25).Value = -1
$ws.Cells.Item(158, 26).Value = 0
$ws.Cells.Item(158, 27).Value = -0
$ws.Cells.Item(158, 28).Value = 0
$ws.Cells.Item(158, 29).Value = -0

# --- Add new rows 159-167 (new matches) ---
# Row 159
$ws.Cells.Item(159, 1).Value = 157
$ws.Cells.Item(159, 2).Value = 7990776
$ws.Cells.Item(159, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(159, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(159, 5).Value = 45375.66666666666
$ws.Cells.Item(159, 6).Value = "Danubio"
$ws.Cells.Item(159, 7).Value = "Montevideo Wanderers"
$ws.Cells.Item(159, 8).Value = 0
$ws.Cells.Item(159, 9).Value = 2
$ws.Cells.Item(159, 10).Value = "A"
$ws.Cells.Item(159, 11).Value = 1.833
$ws.Cells.Item(159, 12).Value = 3.4
$ws.Cells.Item(159, 13).Value = 4.333
$ws.Cells.Item(159, 14).Value = 2.05
$ws.Cells.Item(159, 15).Value = 3.25
$ws.Cells.Item(159, 16).Value = 3.6
$ws.Cells.Item(159, 17).Value = -0.25
$ws.Cells.Item(159, 18).Value = 1.8
$ws.Cells.Item(159, 19).Value = 2.05
$ws.Cells.Item(159, 20).Value = 2
$ws.Cells.Item(159, 21).Value = 1.925
$ws.Cells.Item(159, 22).Value = 1.925
$ws.Cells.Item(159, 23).Value = -1
$ws.Cells.Item(159, 24).Value = -1
$ws.Cells.Item(159, 25).Value = 2.6
$ws.Cells.Item(159, 26).Value = -1
$ws.Cells.Item(159, 27).Value = 1.05
$ws.Cells.Item(159, 28).Value = 0
$ws.Cells.Item(159, 29).Value = -0
$ws.Range("A2").Copy()
$ws.Range("A159").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E159").PasteSpecial(-4122)

# Row 160
$ws.Cells.Item(160, 1).Value = 158
$ws.Cells.Item(160, 2).Value = 7994681
$ws.Cells.Item(160, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(160, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(160, 5).Value = 45380.52083333334
$ws.Cells.Item(160, 6).Value = "Cerro Largo"
$ws.Cells.Item(160, 7).Value = "Racing Club de Montevideo"
$ws.Cells.Item(160, 11).Value = 2.3
$ws.Cells.Item(160, 12).Value = 3
$ws.Cells.Item(160, 13).Value = 3.1
$ws.Cells.Item(160, 14).Value = 2.375
$ws.Cells.Item(160, 15).Value = 3
$ws.Cells.Item(160, 16).Value = 3
$ws.Cells.Item(160, 17).Value = -0.25
$ws.Cells.Item(160, 18).Value = 2.1
$ws.Cells.Item(160, 19).Value = 1.775
$ws.Cells.Item(160, 20).Value = 2.25
$ws.Cells.Item(160, 21).Value = 2.05
$ws.Cells.Item(160, 22).Value = 1.8
$ws.Cells.Item(160, 23).Value = 0
$ws.Cells.Item(160, 24).Value = 0
$ws.Cells.Item(160, 25).Value = 0
$ws.Cells.Item(160, 26).Value = 0
$ws.Cells.Item(160, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Range("A160").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E160").PasteSpecial(-4122)

# Row 161
$ws.Cells.Item(161, 1).Value = 159
$ws.Cells.Item(161, 2).Value = 7994680
$ws.Cells.Item(161, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(161, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(161, 5).Value = 45380.66666666666
$ws.Cells.Item(161, 6).Value = "Penarol"
$ws.Cells.Item(161, 7).Value = "Nacional De Football"
$ws.Cells.Item(161, 11).Value = 2.25
$ws.Cells.Item(161, 12).Value = 3.1
$ws.Cells.Item(161, 13).Value = 3
$ws.Cells.Item(161, 14).Value = 2.15
$ws.Cells.Item(161, 15).Value = 3.1
$ws.Cells.Item(161, 16).Value = 3.2
$ws.Cells.Item(161, 17).Value = -0.25
$ws.Cells.Item(161, 18).Value = 1.925
$ws.Cells.Item(161, 19).Value = 1.925
$ws.Cells.Item(161, 20).Value = 2.25
$ws.Cells.Item(161, 21).Value = 1.875
$ws.Cells.Item(161, 22).Value = 1.975
$ws.Cells.Item(161, 23).Value = 0
$ws.Cells.Item(161, 24).Value = 0
$ws.Cells.Item(161, 25).Value = 0
$ws.Cells.Item(161, 26).Value = 0
$ws.Cells.Item(161, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Range("A161").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E161").PasteSpecial(-4122)

# Row 162
$ws.Cells.Item(162, 1).Value = 160
$ws.Cells.Item(162, 2).Value = 7994520
$ws.Cells.Item(162, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(162, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(162, 5).Value = 45381.5625
$ws.Cells.Item(162, 6).Value = "Rampla Juniors"
$ws.Cells.Item(162, 7).Value = "Atletico Fenix Montevideo"
$ws.Cells.Item(162, 11).Value = 2.8
$ws.Cells.Item(162, 12).Value = 3
$ws.Cells.Item(162, 13).Value = 2.6
$ws.Cells.Item(162, 14).Value = 3
$ws.Cells.Item(162, 15).Value = 3
$ws.Cells.Item(162, 16).Value = 2.45
$ws.Cells.Item(162, 17).Value = 0
$ws.Cells.Item(162, 18).Value = 2.125
$ws.Cells.Item(162, 19).Value = 1.75
$ws.Cells.Item(162, 20).Value = 2.25
$ws.Cells.Item(162, 21).Value = 1.95
$ws.Cells.Item(162, 22).Value = 1.9
$ws.Cells.Item(162, 23).Value = 0
$ws.Cells.Item(162, 24).Value = 0
$ws.Cells.Item(162, 25).Value = 0
$ws.Cells.Item(162, 26).Value = 0
$ws.Cells.Item(162, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Range("A162").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E162").PasteSpecial(-4122)

# Row 163
$ws.Cells.Item(163, 1).Value = 161
$ws.Cells.Item(163, 2).Value = 7994682
$ws.Cells.Item(163, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(163, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(163, 5).Value = 45381.66666666666
$ws.Cells.Item(163, 6).Value = "Liverpool Montevideo"
$ws.Cells.Item(163, 7).Value = "CA River Plate"
$ws.Cells.Item(163, 11).Value = 1.909
$ws.Cells.Item(163, 12).Value = 3.2
$ws.Cells.Item(163, 13).Value = 4
$ws.Cells.Item(163, 14).Value = 1.95
$ws.Cells.Item(163, 15).Value = 3.1
$ws.Cells.Item(163, 16).Value = 4
$ws.Cells.Item(163, 17).Value = -0.5
$ws.Cells.Item(163, 18).Value = 2
$ws.Cells.Item(163, 19).Value = 1.85
$ws.Cells.Item(163, 20).Value = 2.25
$ws.Cells.Item(163, 21).Value = 1.925
$ws.Cells.Item(163, 22).Value = 1.925
$ws.Cells.Item(163, 23).Value = 0
$ws.Cells.Item(163, 24).Value = 0
$ws.Cells.Item(163, 25).Value = 0
$ws.Cells.Item(163, 26).Value = 0
$ws.Cells.Item(163, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Range("A163").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E163").PasteSpecial(-4122)

# Row 164
$ws.Cells.Item(164, 1).Value = 162
$ws.Cells.Item(164, 2).Value = 7994684
$ws.Cells.Item(164, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(164, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(164, 5).Value = 45381.79166666666
$ws.Cells.Item(164, 6).Value = "Defensor Sporting"
$ws.Cells.Item(164, 7).Value = "Danubio"
$ws.Cells.Item(164, 11).Value = 1.909
$ws.Cells.Item(164, 12).Value = 3.25
$ws.Cells.Item(164, 13).Value = 3.8
$ws.Cells.Item(164, 14).Value = 1.909
$ws.Cells.Item(164, 15).Value = 3.25
$ws.Cells.Item(164, 16).Value = 3.8
$ws.Cells.Item(164, 17).Value = -0.5
$ws.Cells.Item(164, 18).Value = 1.925
$ws.Cells.Item(164, 19).Value = 1.925
$ws.Cells.Item(164, 20).Value = 2.25
$ws.Cells.Item(164, 21).Value = 1.9
$ws.Cells.Item(164, 22).Value = 1.95
$ws.Cells.Item(164, 23).Value = 0
$ws.Cells.Item(164, 24).Value = 0
$ws.Cells.Item(164, 25).Value = 0
$ws.Cells.Item(164, 26).Value = 0
$ws.Cells.Item(164, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Range("A164").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E164").PasteSpecial(-4122)

# Row 165
$ws.Cells.Item(165, 1).Value = 163
$ws.Cells.Item(165, 2).Value = 7995146
$ws.Cells.Item(165, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(165, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(165, 5).Value = 45382.625
$ws.Cells.Item(165, 6).Value = "Club Atletico Progreso"
$ws.Cells.Item(165, 7).Value = "Deportivo Maldonado"
$ws.Cells.Item(165, 11).Value = 2.4
$ws.Cells.Item(165, 12).Value = 3.1
$ws.Cells.Item(165, 13).Value = 3
$ws.Cells.Item(165, 14).Value = 2.4
$ws.Cells.Item(165, 15).Value = 3.2
$ws.Cells.Item(165, 16).Value = 2.9
$ws.Cells.Item(165, 17).Value = -0.25
$ws.Cells.Item(165, 18).Value = 2.1
$ws.Cells.Item(165, 19).Value = 1.775
$ws.Cells.Item(165, 20).Value = 2.25
$ws.Cells.Item(165, 21).Value = 1.8
$ws.Cells.Item(165, 22).Value = 2.05
$ws.Cells.Item(165, 23).Value = 0
$ws.Cells.Item(165, 24).Value = 0
$ws.Cells.Item(165, 25).Value = 0
$ws.Cells.Item(165, 26).Value = 0
$ws.Cells.Item(165, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Range("A165").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E165").PasteSpecial(-4122)

# Row 166
$ws.Cells.Item(166, 1).Value = 164
$ws.Cells.Item(166, 2).Value = 7995141
$ws.Cells.Item(166, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(166, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(166, 5).Value = 45382.72916666666
$ws.Cells.Item(166, 6).Value = "Miramar Misiones"
$ws.Cells.Item(166, 7).Value = "Cerro"
$ws.Cells.Item(166, 11).Value = 2.6
$ws.Cells.Item(166, 12).Value = 3
$ws.Cells.Item(166, 13).Value = 2.75
$ws.Cells.Item(166, 14).Value = 2.5
$ws.Cells.Item(166, 15).Value = 3
$ws.Cells.Item(166, 16).Value = 2.875
$ws.Cells.Item(166, 17).Value = 0
$ws.Cells.Item(166, 18).Value = 1.8
$ws.Cells.Item(166, 19).Value = 2.05
$ws.Cells.Item(166, 20).Value = 2.25
$ws.Cells.Item(166, 21).Value = 1.975
$ws.Cells.Item(166, 22).Value = 1.875
$ws.Cells.Item(166, 23).Value = 0
$ws.Cells.Item(166, 24).Value = 0
$ws.Cells.Item(166, 25).Value = 0
$ws.Cells.Item(166, 26).Value = 0
$ws.Cells.Item(166, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Range("A166").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E166").PasteSpecial(-4122)

# Row 167
$ws.Cells.Item(167, 1).Value = 165
$ws.Cells.Item(167, 2).Value = 7994683
$ws.Cells.Item(167, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(167, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(167, 5).Value = 45382.83333333334
$ws.Cells.Item(167, 6).Value = "Montevideo Wanderers"
$ws.Cells.Item(167, 7).Value = "Boston River"
$ws.Cells.Item(167, 11).Value = 2.5
$ws.Cells.Item(167, 12).Value = 3.1
$ws.Cells.Item(167, 13).Value = 2.75
$ws.Cells.Item(167, 14).Value = 2.7
$ws.Cells.Item(167, 15).Value = 3
$ws.Cells.Item(167, 16).Value = 2.6
$ws.Cells.Item(167, 17).Value = 0
$ws.Cells.Item(167, 18).Value = 1.975
$ws.Cells.Item(167, 19).Value = 1.875
$ws.Cells.Item(167, 20).Value = 2.25
$ws.Cells.Item(167, 21).Value = 2
$ws.Cells.Item(167, 22).Value = 1.85
$ws.Cells.Item(167, 23).Value = 0
$ws.Cells.Item(167, 24).Value = 0
$ws.Cells.Item(167, 25).Value = 0
$ws.Cells.Item(167, 26).Value = 0
$ws.Cells.Item(167, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Range("A167").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E167").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Output "Edit complete"